# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps to reflect a new handback report run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row
$overview.Range("G2").Value = "2016-09-01 19:16:52"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for first file row
$zhcn.Range("H2").Value = "2016-09-01 19:16:47"
$zhcn.Range("K2").Value = "2016-09-01 19:17:14"

# de-de sheet: Correspond Handoff Datetime (shared with Overview G2) / Correspond Handback DateTime
$dede.Range("H2").Value = "2016-09-01 19:16:52"
$dede.Range("K2").Value = "2016-09-01 19:17:22"
